$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "NO"
$ws.Range("L2").Value = "SKIPPED"
$ws.Range("M2").Value = "TO BE EXECUTED = NO"
$ws.Range("N2:P2").ClearContents()

# Row 7
$ws.Range("A7").Value = "NO"
$ws.Range("L7").Value = "SKIPPED"
$ws.Range("M7").Value = "TO BE EXECUTED = NO"
$ws.Range("N7:P7").ClearContents()

# Row 8
$ws.Range("A8").Value = "NO"
$ws.Range("L8").Value = "SKIPPED"
$ws.Range("M8").Value = "TO BE EXECUTED = NO"
$ws.Range("N8:P8").ClearContents()

# Row 9
$ws.Range("A9").Value = "NO"
$ws.Range("L9").Value = "SKIPPED"
$ws.Range("M9").Value = "TO BE EXECUTED = NO"
$ws.Range("N9:P9").ClearContents()

# Row 10
$ws.Range("A10").Value = "NO"
$ws.Range("L10").Value = "SKIPPED"
$ws.Range("M10").Value = "TO BE EXECUTED = NO"
$ws.Range("N10:P10").ClearContents()

# Row 11
$ws.Range("A11").Value = "NO"
$ws.Range("L11").Value = "SKIPPED"
$ws.Range("M11").Value = "TO BE EXECUTED = NO"
$ws.Range("N11:P11").ClearContents()

# Row 12
$ws.Range("A12").Value = "NO"
$ws.Range("L12").Value = "SKIPPED"
$ws.Range("M12").Value = "TO BE EXECUTED = NO"
$ws.Range("N12:P12").ClearContents()

# Row 13
$ws.Range("A13").Value = "NO"
$ws.Range("L13").Value = "SKIPPED"
$ws.Range("M13").Value = "TO BE EXECUTED = NO"
$ws.Range("N13:P13").ClearContents()

# Row 14
$ws.Range("A14").Value = "NO"
$ws.Range("L14").Value = "SKIPPED"
$ws.Range("M14").Value = "TO BE EXECUTED = NO"
$ws.Range("N14:P14").ClearContents()

# Row 15
$ws.Range("A15").Value = "NO"
$ws.Range("L15").Value = "SKIPPED"
$ws.Range("M15").Value = "TO BE EXECUTED = NO"
$ws.Range("N15:P15").ClearContents()

# Row 16
$ws.Range("A16").Value = "NO"
$ws.Range("L16").Value = "SKIPPED"
$ws.Range("M16").Value = "TO BE EXECUTED = NO"
$ws.Range("N16:P16").ClearContents()

# Row 17
$ws.Range("A17").Value = "NO"
$ws.Range("L17").Value = "SKIPPED"
$ws.Range("M17").Value = "TO BE EXECUTED = NO"
$ws.Range("N17:P17").ClearContents()

# Row 18
$ws.Range("A18").Value = "NO"
$ws.Range("L18").Value = "SKIPPED"
$ws.Range("M18").Value = "TO BE EXECUTED = NO"

# Row 19
$ws.Range("A19").Value = "NO"
$ws.Range("L19").Value = "SKIPPED"
$ws.Range("M19").Value = "TO BE EXECUTED = NO"

# Row 20
$ws.Range("A20").Value = "NO"
$ws.Range("L20").Value = "SKIPPED"
$ws.Range("M20").Value = "TO BE EXECUTED = NO"

# Row 21
$ws.Range("A21").Value = "NO"
$ws.Range("L21").Value = "SKIPPED"
$ws.Range("M21").Value = "TO BE EXECUTED = NO"

# Row 22
$ws.Range("A22").Value = "NO"
$ws.Range("L22").Value = "SKIPPED"
$ws.Range("M22").Value = "TO BE EXECUTED = NO"

# Row 23
$ws.Range("A23").Value = "NO"
$ws.Range("L23").Value = "SKIPPED"
$ws.Range("M23").Value = "TO BE EXECUTED = NO"

# Row 24
$ws.Range("A24").Value = "NO"
$ws.Range("L24").Value = "SKIPPED"
$ws.Range("M24").Value = "TO BE EXECUTED = NO"

# Row 25
$ws.Range("L25").Value = "SKIPPED"
$ws.Range("M25").Value = "TO BE EXECUTED = NO"

# Row 26
$ws.Range("L26").Value = "SKIPPED"
$ws.Range("M26").Value = "TO BE EXECUTED = NO"

# Row 27
$ws.Range("L27").Value = "SKIPPED"
$ws.Range("M27").Value = "TO BE EXECUTED = NO"

# Row 28
$ws.Range("L28").Value = "SKIPPED"
$ws.Range("M28").Value = "TO BE EXECUTED = NO"

# Row 29
$ws.Range("L29").Value = "SKIPPED"
$ws.Range("M29").Value = "TO BE EXECUTED = NO"

# Row 30
$ws.Range("L30").Value = "SKIPPED"
$ws.Range("M30").Value = "TO BE EXECUTED = NO"

# Row 31
$ws.Range("L31").Value = "SKIPPED"
$ws.Range("M31").Value = "TO BE EXECUTED = NO"

# Row 32
$ws.Range("L32").Value = "SKIPPED"
$ws.Range("M32").Value = "TO BE EXECUTED = NO"

# Row 33
$ws.Range("L33").Value = "SKIPPED"
$ws.Range("M33").Value = "TO BE EXECUTED = NO"

# Row 34
$ws.Range("A34").Value = "YES"
$ws.Range("L34").Value = "PASS"
$ws.Range("M34").Formula = '=""'
$ws.Range("N34").Value = "Opened: https://www.flipkart.com/"
$ws.Range("O34").Value = "screenshots/STEP_33.png"
$ws.Range("P34").Value = "page_sources/STEP_33_source.html"

# Row 35
$ws.Range("A35").Value = "YES"
$ws.Range("L35").Value = "PASS"
$ws.Range("M35").Formula = '=""'
$ws.Range("N35").Value = "Clicked: Login"
$ws.Range("O35").Value = "screenshots/STEP_34.png"
$ws.Range("P35").Value = "page_sources/STEP_34_source.html"

# Row 36
$ws.Range("A36").Value = "YES"
$ws.Range("L36").Value = "PASS"
$ws.Range("M36").Formula = '=""'
$ws.Range("N36").Value = "Clicked: Men"
$ws.Range("O36").Value = "screenshots/STEP_35.png"
$ws.Range("P36").Value = "page_sources/STEP_35_source.html"

# Row 37
$ws.Range("A37").Value = "YES"
$ws.Range("L37").Value = "PASS"
$ws.Range("M37").Formula = '=""'
$ws.Range("N37").Value = "Clicked: Running Shoes"
$ws.Range("O37").Value = "screenshots/STEP_36.png"
$ws.Range("P37").Value = "page_sources/STEP_36_source.html"

# Row 38
$ws.Range("A38").Value = "YES"
$ws.Range("L38").Value = "PASS"
$ws.Range("M38").Formula = '=""'
$ws.Range("N38").Value = "Clicked: COLOR"
$ws.Range("O38").Value = "screenshots/STEP_37.png"
$ws.Range("P38").Value = "page_sources/STEP_37_source.html"

# Row 39
$ws.Range("A39").Value = "YES"
$ws.Range("L39").Value = "FAIL"
$ws.Range("M39").Value = "Could not click element"
$ws.Range("N39").Value = "Failed to click: Blue"
$ws.Range("O39").Value = "screenshots/STEP_38.png"
$ws.Range("P39").Value = "page_sources/STEP_38_source.html"

# Row 40
$ws.Range("A40").Value = "YES"
$ws.Range("L40").Value = "FAIL"
$ws.Range("M40").Value = "Could not click element"
$ws.Range("N40").Value = "Failed to click: BestSeller"
$ws.Range("O40").Value = "screenshots/STEP_39.png"
$ws.Range("P40").Value = "page_sources/STEP_39_source.html"
